$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-28 with new ticker values ---
$ws.Range("B2").Value = "NSE:AJMERA"
$ws.Range("C2").Value = "NSE:ADSL"
$ws.Range("D2").Value = "NSE:BHEL"
$ws.Range("F2").Value = "NSE:BEL"
$ws.Range("B3").Value = "NSE:ALKEM"
$ws.Range("C3").Value = "NSE:ATL"
$ws.Range("D3").Value = "NSE:NATIONALUM"
$ws.Range("F3").Value = "NSE:IRCTC"
$ws.Range("B4").Value = "NSE:ALPHAETF"
$ws.Range("C4").Value = "NSE:CHEMPLASTS"
$ws.Range("F4").Value = "NSE:ITC"
$ws.Range("B5").Value = "NSE:ANANDRATHI"
$ws.Range("C5").Value = "NSE:DELHIVERY"
$ws.Range("F5").Value = "NSE:NATIONALUM"
$ws.Range("B6").Value = "NSE:ANANTRAJ"
$ws.Range("C6").Value = "NSE:DLINKINDIA"
$ws.Range("B7").Value = "NSE:BEL"
$ws.Range("C7").Value = "NSE:ENIL"
$ws.Range("B8").Value = "NSE:GESHIP"
$ws.Range("C8").Value = "NSE:GPPL"
$ws.Range("B9").Value = "NSE:HBLPOWER"
$ws.Range("C9").Value = "NSE:JKPAPER"
$ws.Range("B10").Value = "NSE:HDFCBSE500"
$ws.Range("C10").Value = "NSE:KAKATCEM"
$ws.Range("B11").Value = "NSE:HDFCSML250"
$ws.Range("C11").Value = "NSE:KINGFA"
$ws.Range("B12").Value = "NSE:HINDCOPPER"
$ws.Range("C12").Value = "NSE:LOVABLE"
$ws.Range("B13").Value = "NSE:HPAL"
$ws.Range("C13").Value = "NSE:LTIM"
$ws.Range("B14").Value = "NSE:IRCTC"
$ws.Range("C14").Value = "NSE:MARALOVER"
$ws.Range("B15").Value = "NSE:ITC"
$ws.Range("C15").Value = "NSE:NUVAMA"
$ws.Range("B16").Value = "NSE:ITI"
$ws.Range("C16").Value = "NSE:POKARNA"
$ws.Range("B17").Value = "NSE:IVZINGOLD"
$ws.Range("C17").ClearContents()
$ws.Range("B18").Value = "NSE:JBMA"
$ws.Range("B19").Value = "NSE:JITFINFRA"
$ws.Range("B20").Value = "NSE:KELLTONTEC"
$ws.Range("B21").Value = "NSE:KIRIINDUS"
$ws.Range("B22").Value = "NSE:KSCL"
$ws.Range("B23").Value = "NSE:LEXUS"
$ws.Range("B24").Value = "NSE:LPDC"
$ws.Range("B25").Value = "NSE:MIDHANI"
$ws.Range("B26").Value = "NSE:MIRCELECTR"
$ws.Range("B27").Value = "NSE:MON100"
$ws.Range("B28").Value = "NSE:MONQ50"

# --- Append new rows 29-34 (copy column-A formatting from row 28) ---
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "NSE:NATIONALUM"
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "NSE:NETWEB"
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "NSE:NV20BEES"
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "NSE:PAGEIND"
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "NSE:PNC"
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "NSE:RADIANTCMS"

Write-Output "edit applied"
